# "Test no volume ksvmeans"
#
# - Adds a "No Volume" results column (F) to the K-SVMeans sheet, with an
#   AVERAGE formula summarizing it (matching the existing B/C/E columns).
# - Moves the active/selected tab from "ANN" (sheet 1) to "K-SVMeans" (sheet 3).
# - Updates the remembered cell-selection on each sheet.
# - Turns on explicit portrait page setup for the K-SVMeans sheet.

$wb = $excel.ActiveWorkbook

$wsAnn = $wb.Worksheets.Item(1)      # "ANN"
$wsSvm = $wb.Worksheets.Item(2)      # "SVM"
$wsKsm = $wb.Worksheets.Item(3)      # "K-SVMeans"

# --- New "No Volume" (F) column data on the K-SVMeans sheet ---------------
$wsKsm.Range("F4").Value = 59.42
$wsKsm.Range("F5").Value = 55.67

$wsKsm.Range("F6").Value = 51.87
# F6 picks up the same "muted" style already used by E8 in this table.
$wsKsm.Range("E8").Copy()
$wsKsm.Range("F6").PasteSpecial(-4122)  # xlPasteFormats

$wsKsm.Range("F7").Value = 57.61
$wsKsm.Range("F8").Value = 55.84

# Average row, styled like the existing B9/C9/E9 total cells.
$wsKsm.Range("F9").Formula = "=AVERAGE(F4:F8)"
$wsKsm.Range("E9").Copy()
$wsKsm.Range("F9").PasteSpecial(-4122)  # xlPasteFormats

$wsKsm.Application.CutCopyMode = $false

# --- Explicit page setup on the K-SVMeans sheet ----------------------------
$wsKsm.PageSetup.Orientation = 1  # xlPortrait

# --- Selection / active-sheet bookkeeping ----------------------------------
# Update the remembered selection on ANN and SVM without leaving them as the
# active tab.
$wsAnn.Activate()
$wsAnn.Range("H9").Select()

$wsSvm.Activate()
$wsSvm.Range("G6").Select()

# K-SVMeans becomes the active (selected) tab, with its own new selection.
$wsKsm.Activate()
$wsKsm.Range("H11").Select()
